$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45033
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 861
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 45028
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("R3").Value = 'Región Metropolitana'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 45020
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = '$/caja 16 kilos'
$ws.Range("R4").Value = 'Provincia de Los Andes'
$ws.Range("S4").Value = 938
$ws.Range("T4").Value = 16

# Row 5
$ws.Range("D5").Value = 45002
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 45001
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 972
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44999
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 17500
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 972
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 45037
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 16000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 889
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 45014
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 18000
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 45021
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("R10").Value = 'Provincia de Los Andes'
$ws.Range("S10").Value = 861
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 45030
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 18000
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 45041
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 833
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 45036
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15500
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 861
$ws.Range("T13").Value = 18

